# Auto-generated script applying scheduled market-data refresh to Sheets workbook
# Updates currentAveragePrice*/LevePrice*/LeveProfit* columns (H:N) on each sheet
# to reflect the latest pulled market data, per commit "chore: update Sheets via scheduled runner".

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1232.7887
$ws.Range("I15").Value = 1232.7887
$ws.Range("K15").Value = 3698.3661
$ws.Range("M15").Value = -3529.3661
$ws.Range("H17").Value = 1591.909
$ws.Range("J17").Value = 1762.3043
$ws.Range("L17").Value = 5286.9129
$ws.Range("N17").Value = -5622.9129
$ws.Range("H32").Value = 3678.4375
$ws.Range("I32").Value = 3777.5557
$ws.Range("J32").Value = 3551
$ws.Range("K32").Value = 3777.5557
$ws.Range("L32").Value = 3551
$ws.Range("M32").Value = -3451.5557
$ws.Range("N32").Value = -4203
$ws.Range("H40").Value = 5428.2144
$ws.Range("I40").Value = 8665.833000000001
$ws.Range("K40").Value = 8665.833000000001
$ws.Range("M40").Value = -8490.833000000001
$ws.Range("H61").Value = 1507.5
$ws.Range("I61").Value = 1507.5
$ws.Range("K61").Value = 4522.5
$ws.Range("M61").Value = -4350.5
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()  # was -10494.667
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()  # was -11714.667
$ws.Range("H74").Value = 7490.3335
$ws.Range("I74").Value = 5842
$ws.Range("K74").Value = 5842
$ws.Range("M74").Value = -4906
$ws.Range("H77").Value = 7490.3335
$ws.Range("I77").Value = 5842
$ws.Range("K77").Value = 29210
$ws.Range("M77").Value = -24530
$ws.Range("H86").Value = 2674.5881
$ws.Range("I86").Value = 2683.6428
$ws.Range("K86").Value = 2683.6428
$ws.Range("M86").Value = -1560.6428
$ws.Range("H89").Value = 2674.5881
$ws.Range("I89").Value = 2683.6428
$ws.Range("K89").Value = 13418.214
$ws.Range("M89").Value = -7802.214
$ws.Range("H92").Value = 1093.2222
$ws.Range("I92").Value = 338.2
$ws.Range("J92").Value = 4868.3335
$ws.Range("K92").Value = 338.2
$ws.Range("L92").Value = 4868.3335
$ws.Range("M92").Value = 909.8
$ws.Range("N92").Value = -7364.3335
$ws.Range("H94").Value = 3543.1538
$ws.Range("I94").Value = 1732.8182
$ws.Range("K94").Value = 1732.8182
$ws.Range("M94").Value = -1281.8182
$ws.Range("H96").Value = 531.5714
$ws.Range("I96").Value = 536.8333
$ws.Range("J96").Value = 500
$ws.Range("K96").Value = 1610.4999
$ws.Range("L96").Value = 1500
$ws.Range("M96").Value = -237.4999
$ws.Range("N96").Value = -4246
$ws.Range("H97").Value = 1505.2222
$ws.Range("J97").Value = 1655.875
$ws.Range("L97").Value = 4967.625
$ws.Range("N97").Value = -5959.625
$ws.Range("H104").Value = 110
$ws.Range("I104").Value = 110
$ws.Range("K104").Value = 330
$ws.Range("M104").Value = 1417
$ws.Range("H106").Value = 1863.2162
$ws.Range("I106").Value = 1176.3928
$ws.Range("K106").Value = 1176.3928
$ws.Range("M106").Value = -545.3928000000001
$ws.Range("H107").Value = 445.22223
$ws.Range("J107").Value = 452.2
$ws.Range("L107").Value = 452.2
$ws.Range("N107").Value = -4292.2
$ws.Range("H112").Value = 1934.381
$ws.Range("J112").Value = 2006.4736
$ws.Range("L112").Value = 6019.4208
$ws.Range("N112").Value = -8235.4208
$ws.Range("H113").Value = 5740.1177
$ws.Range("I113").Value = 4222.3335
$ws.Range("J113").Value = 7447.625
$ws.Range("K113").Value = 4222.3335
$ws.Range("L113").Value = 7447.625
$ws.Range("M113").Value = -968.3334999999997
$ws.Range("N113").Value = -13955.625
$ws.Range("H131").Value = 911.25
$ws.Range("I131").Value = 970
$ws.Range("J131").Value = 500
$ws.Range("K131").Value = 2910
$ws.Range("L131").Value = 1500
$ws.Range("M131").Value = 2130
$ws.Range("N131").Value = -11580
$ws.Range("H132").Value = 384317.3
$ws.Range("I132").Value = 409816.1
$ws.Range("K132").Value = 1229448.3
$ws.Range("M132").Value = -1226918.3
$ws.Range("H135").Value = 41667424
$ws.Range("I135").Value = 45455280
$ws.Range("K135").Value = 409097520
$ws.Range("M135").Value = -409094985
$ws.Range("H137").Value = 3966.2856
$ws.Range("I137").Value = 2899.6
$ws.Range("K137").Value = 8698.799999999999
$ws.Range("M137").Value = -6148.799999999999
$ws.Range("H138").Value = 4145.457
$ws.Range("I138").Value = 2317.25
$ws.Range("J138").Value = 5685
$ws.Range("K138").Value = 6951.75
$ws.Range("L138").Value = 17055
$ws.Range("M138").Value = -1811.75
$ws.Range("N138").Value = -27335

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1840586.5
$ws.Range("I32").Value = 1967486.6
$ws.Range("K32").Value = 1967486.6
$ws.Range("M32").Value = -1967199.6
$ws.Range("H36").Value = 20000
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 20000
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 20000
$ws.Range("M36").ClearContents()  # was -4654
$ws.Range("N36").Value = -20692
$ws.Range("H45").Value = 2797.6086
$ws.Range("I45").Value = 2644.9
$ws.Range("J45").Value = 3815.6667
$ws.Range("K45").Value = 2644.9
$ws.Range("L45").Value = 3815.6667
$ws.Range("M45").Value = -2267.9
$ws.Range("N45").Value = -4569.6667
$ws.Range("H74").Value = 2944120.5
$ws.Range("I74").Value = 4465851.5
$ws.Range("J74").Value = 5605
$ws.Range("K74").Value = 4465851.5
$ws.Range("L74").Value = 5605
$ws.Range("M74").Value = -4464977.5
$ws.Range("N74").Value = -7353
$ws.Range("H77").Value = 2944120.5
$ws.Range("I77").Value = 4465851.5
$ws.Range("J77").Value = 5605
$ws.Range("K77").Value = 22329257.5
$ws.Range("L77").Value = 28025
$ws.Range("M77").Value = -22324889.5
$ws.Range("N77").Value = -36761
$ws.Range("H102").Value = 3313.4375
$ws.Range("I102").Value = 3313.4375
$ws.Range("K102").Value = 3313.4375
$ws.Range("M102").Value = -1691.4375
$ws.Range("H110").Value = 1884.4642
$ws.Range("I110").Value = 1925.1154
$ws.Range("J110").Value = 1356
$ws.Range("K110").Value = 1925.1154
$ws.Range("L110").Value = 1356
$ws.Range("M110").Value = 119.8846000000001
$ws.Range("N110").Value = -5446
$ws.Range("H122").Value = 2545.9429
$ws.Range("I122").Value = 2565.9062
$ws.Range("K122").Value = 7697.7186
$ws.Range("M122").Value = -5247.7186
$ws.Range("H132").Value = 287027.44
$ws.Range("I132").Value = 486177.22
$ws.Range("J132").Value = 5874.853
$ws.Range("K132").Value = 1458531.66
$ws.Range("L132").Value = 17624.559
$ws.Range("M132").Value = -1456001.66
$ws.Range("N132").Value = -22684.559
$ws.Range("H135").Value = 90000
$ws.Range("J135").Value = 90000
$ws.Range("L135").Value = 90000
$ws.Range("N135").Value = -100140

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1493.7693
$ws.Range("I86").Value = 1572
$ws.Range("J86").Value = 1233
$ws.Range("K86").Value = 1572
$ws.Range("L86").Value = 1233
$ws.Range("M86").Value = -449
$ws.Range("N86").Value = -3479
$ws.Range("H89").Value = 1493.7693
$ws.Range("I89").Value = 1572
$ws.Range("J89").Value = 1233
$ws.Range("K89").Value = 7860
$ws.Range("L89").Value = 6165
$ws.Range("M89").Value = -2244
$ws.Range("N89").Value = -17397
$ws.Range("H99").Value = 2515.077
$ws.Range("I99").Value = 1516.0555
$ws.Range("J99").Value = 4762.875
$ws.Range("K99").Value = 1516.0555
$ws.Range("L99").Value = 4762.875
$ws.Range("M99").Value = -18.05549999999994
$ws.Range("N99").Value = -7758.875
$ws.Range("H105").Value = 1999
$ws.Range("I105").Value = 1624
$ws.Range("J105").Value = 2499
$ws.Range("K105").Value = 1624
$ws.Range("L105").Value = 2499
$ws.Range("M105").Value = 123
$ws.Range("N105").Value = -5993
$ws.Range("H107").Value = 779.4
$ws.Range("I107").Value = 779.4
$ws.Range("K107").Value = 779.4
$ws.Range("M107").Value = 1140.6
$ws.Range("H134").Value = 598126
$ws.Range("I134").Value = 852255.25
$ws.Range("K134").Value = 2556765.75
$ws.Range("M134").Value = -2554230.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 15874.5
$ws.Range("J11").Value = 15874.5
$ws.Range("L11").Value = 15874.5
$ws.Range("N11").Value = -16154.5
$ws.Range("H16").Value = 21199.8
$ws.Range("I16").Value = 21199.8
$ws.Range("K16").Value = 21199.8
$ws.Range("M16").Value = -20912.8
$ws.Range("H31").Value = 5397.793
$ws.Range("I31").Value = 2863.389
$ws.Range("J31").Value = 9545
$ws.Range("K31").Value = 2863.389
$ws.Range("L31").Value = 9545
$ws.Range("M31").Value = -2568.389
$ws.Range("N31").Value = -10135
$ws.Range("H34").Value = 5397.793
$ws.Range("I34").Value = 2863.389
$ws.Range("J34").Value = 9545
$ws.Range("K34").Value = 2863.389
$ws.Range("L34").Value = 9545
$ws.Range("M34").Value = -2661.389
$ws.Range("N34").Value = -9949
$ws.Range("H58").Value = 461149.56
$ws.Range("I58").Value = 540393
$ws.Range("J58").Value = 5499.75
$ws.Range("K58").Value = 540393
$ws.Range("L58").Value = 5499.75
$ws.Range("M58").Value = -540190
$ws.Range("N58").Value = -5905.75
$ws.Range("H69").Value = 26475
$ws.Range("I69").Value = 4219
$ws.Range("J69").Value = 115499
$ws.Range("K69").Value = 4219
$ws.Range("L69").Value = 115499
$ws.Range("M69").Value = -3470
$ws.Range("N69").Value = -116997
$ws.Range("H70").Value = 75788.336
$ws.Range("J70").Value = 75788.336
$ws.Range("L70").Value = 75788.336
$ws.Range("N70").Value = -76418.336
$ws.Range("H72").Value = 26475
$ws.Range("I72").Value = 4219
$ws.Range("J72").Value = 115499
$ws.Range("K72").Value = 12657
$ws.Range("L72").Value = 346497
$ws.Range("M72").Value = -8913
$ws.Range("N72").Value = -353985
$ws.Range("H73").Value = 75788.336
$ws.Range("J73").Value = 75788.336
$ws.Range("L73").Value = 75788.336
$ws.Range("N73").Value = -77972.336
$ws.Range("H86").Value = 2728.1765
$ws.Range("J86").Value = 2799.4285
$ws.Range("L86").Value = 2799.4285
$ws.Range("N86").Value = -5045.4285
$ws.Range("H89").Value = 2728.1765
$ws.Range("J89").Value = 2799.4285
$ws.Range("L89").Value = 13997.1425
$ws.Range("N89").Value = -25229.1425
$ws.Range("H107").Value = 1221.9048
$ws.Range("I107").Value = 877.1053000000001
$ws.Range("K107").Value = 877.1053000000001
$ws.Range("M107").Value = 1042.8947
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()  # was -77180
$ws.Range("H113").Value = 21199.8
$ws.Range("I113").Value = 21199.8
$ws.Range("K113").Value = 21199.8
$ws.Range("M113").Value = -19029.8
$ws.Range("H122").Value = 2553.3635
$ws.Range("I122").Value = 2556.2856
$ws.Range("J122").Value = 2548.25
$ws.Range("K122").Value = 7668.8568
$ws.Range("L122").Value = 7644.75
$ws.Range("M122").Value = -5218.8568
$ws.Range("N122").Value = -12544.75
$ws.Range("H132").Value = 8079614
$ws.Range("I132").Value = 19781.895
$ws.Range("J132").Value = 20841014
$ws.Range("K132").Value = 59345.685
$ws.Range("L132").Value = 62523042
$ws.Range("M132").Value = -56815.685
$ws.Range("N132").Value = -62528102
$ws.Range("H136").Value = 461149.56
$ws.Range("I136").Value = 540393
$ws.Range("J136").Value = 5499.75
$ws.Range("K136").Value = 1621179
$ws.Range("L136").Value = 16499.25
$ws.Range("M136").Value = -1618629
$ws.Range("N136").Value = -21599.25
$ws.Range("H140").Value = 100000
$ws.Range("J140").Value = 100000
$ws.Range("L140").Value = 100000
$ws.Range("N140").Value = -110360

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3188.647
$ws.Range("I5").Value = 673.5454999999999
$ws.Range("K5").Value = 2020.6365
$ws.Range("M5").Value = -1908.6365
$ws.Range("H9").Value = 642.5
$ws.Range("I9").Value = 1000
$ws.Range("K9").Value = 3000
$ws.Range("M9").Value = -2776
$ws.Range("H11").Value = 5264049.5
$ws.Range("I11").Value = 328
$ws.Range("K11").Value = 984
$ws.Range("M11").Value = -844
$ws.Range("H80").Value = 4573.25
$ws.Range("J80").Value = 4559.2
$ws.Range("L80").Value = 13677.6
$ws.Range("N80").Value = -15549.6
$ws.Range("H83").Value = 4573.25
$ws.Range("J83").Value = 4559.2
$ws.Range("L83").Value = 41032.8
$ws.Range("N83").Value = -50392.8
$ws.Range("H113").Value = 3585.7
$ws.Range("I113").Value = 1733
$ws.Range("J113").Value = 4379.7144
$ws.Range("K113").Value = 5199
$ws.Range("L113").Value = 13139.1432
$ws.Range("M113").Value = -3029
$ws.Range("N113").Value = -17479.1432
$ws.Range("H131").Value = 21964.727
$ws.Range("J131").Value = 23998.2
$ws.Range("L131").Value = 71994.60000000001
$ws.Range("N131").Value = -82074.60000000001
$ws.Range("H135").Value = 3188.647
$ws.Range("I135").Value = 673.5454999999999
$ws.Range("K135").Value = 6061.9095
$ws.Range("M135").Value = -3526.9095
$ws.Range("H138").Value = 1936.2858
$ws.Range("I138").Value = 1592.3334
$ws.Range("K138").Value = 4777.0002
$ws.Range("M138").Value = 362.9997999999996

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 267200.16
$ws.Range("J80").Value = 4337.5557
$ws.Range("L80").Value = 4337.5557
$ws.Range("N80").Value = -6333.5557
$ws.Range("H83").Value = 267200.16
$ws.Range("J83").Value = 4337.5557
$ws.Range("L83").Value = 21687.7785
$ws.Range("N83").Value = -31671.7785
$ws.Range("H97").Value = 1056.5807
$ws.Range("I97").Value = 577.7826
$ws.Range("J97").Value = 2433.125
$ws.Range("K97").Value = 577.7826
$ws.Range("L97").Value = 2433.125
$ws.Range("M97").Value = -81.7826
$ws.Range("N97").Value = -3425.125
$ws.Range("H102").Value = 2121.2104
$ws.Range("I102").Value = 1788.4117
$ws.Range("K102").Value = 1788.4117
$ws.Range("M102").Value = -166.4117000000001
$ws.Range("H122").Value = 6335.3794
$ws.Range("I122").Value = 4002
$ws.Range("K122").Value = 12006
$ws.Range("M122").Value = -9556
$ws.Range("H126").Value = 697181.2
$ws.Range("I126").Value = 879613.3
$ws.Range("J126").Value = 3939
$ws.Range("K126").Value = 2638839.9
$ws.Range("L126").Value = 11817
$ws.Range("M126").Value = -2636369.9
$ws.Range("N126").Value = -16757
$ws.Range("H132").Value = 2331.75
$ws.Range("I132").Value = 2288.0386
$ws.Range("J132").Value = 2900
$ws.Range("K132").Value = 6864.1158
$ws.Range("L132").Value = 8700
$ws.Range("M132").Value = -4334.1158
$ws.Range("N132").Value = -13760

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5440.091
$ws.Range("I7").Value = 4878.263
$ws.Range("K7").Value = 4878.263
$ws.Range("M7").Value = -4766.263
$ws.Range("H22").Value = 55815.05
$ws.Range("I22").Value = 251124.75
$ws.Range("K22").Value = 251124.75
$ws.Range("M22").Value = -250829.75
$ws.Range("H27").Value = 55815.05
$ws.Range("I27").Value = 251124.75
$ws.Range("K27").Value = 251124.75
$ws.Range("M27").Value = -251017.75
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()  # was -15820
$ws.Range("H46").Value = 5000
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 5000
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 5000
$ws.Range("M46").ClearContents()  # was -1034
$ws.Range("N46").Value = -5376
$ws.Range("H93").Value = 1911.5
$ws.Range("I93").Value = 1679.4445
$ws.Range("J93").Value = 4000
$ws.Range("K93").Value = 1679.4445
$ws.Range("L93").Value = 4000
$ws.Range("M93").Value = -431.4445000000001
$ws.Range("N93").Value = -6496
$ws.Range("H100").Value = 7537.4707
$ws.Range("I100").Value = 1866.5
$ws.Range("K100").Value = 1866.5
$ws.Range("M100").Value = -1325.5
$ws.Range("H122").Value = 40131.586
$ws.Range("I122").Value = 4310.6665
$ws.Range("J122").Value = 134161.5
$ws.Range("K122").Value = 12931.9995
$ws.Range("L122").Value = 402484.5
$ws.Range("M122").Value = -10481.9995
$ws.Range("N122").Value = -407384.5
$ws.Range("H126").Value = 5440.091
$ws.Range("I126").Value = 4878.263
$ws.Range("K126").Value = 14634.789
$ws.Range("M126").Value = -12164.789
$ws.Range("H132").Value = 1118942
$ws.Range("I132").Value = 1332550.5
$ws.Range("J132").Value = 8178
$ws.Range("K132").Value = 3997651.5
$ws.Range("L132").Value = 24534
$ws.Range("M132").Value = -3995121.5
$ws.Range("N132").Value = -29594
$ws.Range("H136").Value = 3079.8
$ws.Range("I136").Value = 3224.75
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 9674.25
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -7124.25
$ws.Range("N136").Value = -12600

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2994.889
$ws.Range("I81").Value = 2994.889
$ws.Range("K81").Value = 5989.778
$ws.Range("M81").Value = -4928.778
$ws.Range("H84").Value = 2994.889
$ws.Range("I84").Value = 2994.889
$ws.Range("K84").Value = 29948.89
$ws.Range("M84").Value = -24644.89
$ws.Range("H96").Value = 1730.6
$ws.Range("I96").Value = 1329.5834
$ws.Range("K96").Value = 1329.5834
$ws.Range("M96").Value = 43.41660000000002
$ws.Range("H100").Value = 1594.7333
$ws.Range("J100").Value = 1223.625
$ws.Range("L100").Value = 2447.25
$ws.Range("N100").Value = -3529.25
$ws.Range("H113").Value = 2001.6818
$ws.Range("I113").Value = 1594.4546
$ws.Range("K113").Value = 4783.3638
$ws.Range("M113").Value = -2613.3638
$ws.Range("H122").Value = 2834.2812
$ws.Range("I122").Value = 2630.6086
$ws.Range("J122").Value = 3354.7778
$ws.Range("K122").Value = 7891.825800000001
$ws.Range("L122").Value = 10064.3334
$ws.Range("M122").Value = -5441.825800000001
$ws.Range("N122").Value = -14964.3334
$ws.Range("H126").Value = 1983.4166
$ws.Range("I126").Value = 1808.9131
$ws.Range("J126").Value = 5997
$ws.Range("K126").Value = 5426.7393
$ws.Range("L126").Value = 17991
$ws.Range("M126").Value = -2956.7393
$ws.Range("N126").Value = -22931
$ws.Range("H132").Value = 24748184
$ws.Range("I132").Value = 2035954.1
$ws.Range("K132").Value = 6107862.300000001
$ws.Range("M132").Value = -6105332.300000001
$ws.Range("H136").Value = 9048092
$ws.Range("I136").Value = 12666016
$ws.Range("K136").Value = 37998048
$ws.Range("M136").Value = -37995498

